$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.406.33"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.656.43"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.23"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.513"
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.885.59"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.661.93"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.15"
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.545"
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "68.16"
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.372.55"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "221.00"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.88"
$ws.Range("E21").Value = "  +2.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.45"
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.53"
$ws.Range("E23").Value = "  +3.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.27"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.65"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.56"
$ws.Range("E26").Value = "  +1.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.91"
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.05"
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.46"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0178"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.548"
$ws.Range("E38").Value = "  +1.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.844"
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.810"
$ws.Range("E41").Value = "  -1.27%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.23"
$ws.Range("E42").Value = "  +4.76%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.43"
$ws.Range("E43").Value = "  +1.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.795.97"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.38"
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "92.01"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +17.85%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0514"
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.69"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0978"
$ws.Range("E51").Value = "  -0.19%  "
